$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.846.24"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.21"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.44"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("E6").Value = "  -0.55%  "

# Row 7
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("E9").Value = "  -1.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +0.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.871.21"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.50"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -0.10%  "

# Row 15
$ws.Range("E15").Value = "  +0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("E16").Value = "  +1.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.841.54"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("E18").Value = "  -1.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.41"
$ws.Range("E19").Value = "  +1.34%  "

# Row 20
$ws.Range("E20").Value = "  -0.11%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.58"
$ws.Range("E22").Value = "  +5.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.36"
$ws.Range("E23").Value = "  -2.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.80"
$ws.Range("E25").Value = "  +1.93%  "

# Row 26
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  -0.47%  "

# Row 28
$ws.Range("E28").Value = "  +1.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.71"
$ws.Range("E29").Value = "  +0.46%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  -0.40%  "

# Row 31
$ws.Range("E31").Value = "  +1.00%  "

# Row 32
$ws.Range("E32").Value = "  +1.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.284.14"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("E35").Value = "  +0.99%  "

# Row 36
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("E37").Value = "  -0.89%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.534"
$ws.Range("E38").Value = "  +0.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  -0.49%  "

# Row 40
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -0.42%  "

# Row 42
$ws.Range("E42").Value = "  -0.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.782.12"
$ws.Range("E43").Value = "  -0.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  -6.34%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.62"
$ws.Range("E45").Value = "  +1.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.10"
$ws.Range("E46").Value = "  -0.98%  "

# Row 47
$ws.Range("E47").Value = "  -0.55%  "

# Row 48
$ws.Range("E48").Value = "  -1.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.59"
$ws.Range("E49").Value = "  -1.01%  "

# Row 50
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
$ws.Range("E51").Value = "  -0.12%  "
